$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 760
$ws1.Range("F3").Value = 14232
$ws1.Range("F4").Value = 14328
$ws1.Range("F6").Value = 1401
$ws1.Range("F7").Value = 5887
$ws1.Range("F8").Value = 985
$ws1.Range("F9").Value = 574
$ws1.Range("F12").Value = 201
$ws1.Range("F13").Value = 1546
$ws1.Range("F14").Value = 438
$ws1.Range("F16").Value = 1200
$ws1.Range("F17").Value = 1835
$ws1.Range("F20").Value = 2290
$ws1.Range("F21").Value = 565
$ws1.Range("F23").Value = 3325
$ws1.Range("F25").Value = 312
$ws1.Range("F26").Value = 2399
$ws1.Range("F27").Value = 594
$ws1.Range("F30").Value = 1790
$ws1.Range("F32").Value = 1395
$ws1.Range("F33").Value = 101
$ws1.Range("F34").Value = 149
$ws1.Range("F35").Value = 4840
$ws1.Range("F36").Value = 4856
$ws1.Range("F37").Value = 306
$ws1.Range("F40").Value = 682
$ws1.Range("F41").Value = 3294
$ws1.Range("F45").Value = 108
$ws1.Range("F47").Value = 4427
$ws1.Range("F48").Value = 586
$ws1.Range("F49").Value = 292

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 17
$ws2.Range("F22").Value = 56
$ws2.Range("F25").Value = 69

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 7572
$ws3.Range("F3").Value = 241
$ws3.Range("F4").Value = 783

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 760
$ws4.Range("F3").Value = 241
$ws4.Range("F4").Value = 783
$ws4.Range("F6").Value = 14232
$ws4.Range("F8").Value = 1401
$ws4.Range("F9").Value = 5887
$ws4.Range("F10").Value = 985
$ws4.Range("F14").Value = 1546
$ws4.Range("F15").Value = 438
$ws4.Range("F16").Value = 1200
$ws4.Range("F17").Value = 1835
$ws4.Range("F20").Value = 565
$ws4.Range("F21").Value = 3325
$ws4.Range("F22").Value = 312
$ws4.Range("F23").Value = 594
$ws4.Range("F25").Value = 1790
$ws4.Range("F28").Value = 1395
$ws4.Range("F30").Value = 101
$ws4.Range("F31").Value = 149
$ws4.Range("F32").Value = 4840
$ws4.Range("F33").Value = 4856
$ws4.Range("F34").Value = 56
$ws4.Range("F35").Value = 306
$ws4.Range("F38").Value = 682
$ws4.Range("F39").Value = 3294
$ws4.Range("F45").Value = 4427
$ws4.Range("F46").Value = 586
$ws4.Range("F47").Value = 292
